# "Fixed en FrontEnd, Registro de hora de trabajo."
# Adds a new work-log row (row 42) to Hoja1: Bruno Díaz worked 1 hour on
# 2017-05-14 (serial 42869) on "Sprint 3 - FrontEnd", detail "Borrado de
# componentes no no utilizados."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row right after the existing last row (row 41).
$ws.Range("A42").Value = "Bruno Díaz"
$ws.Range("B42").Value = 42869
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = "Sprint 3 - FrontEnd"
$ws.Range("E42").Value = "Borrado de componentes no no utilizados."

# Match the date formatting used by the rest of column B (short date)
# by copying the format from the cell right above instead of assigning a
# brand-new number format (keeps the existing style table intact).
$ws.Range("B41").Copy()
$ws.Range("B42").PasteSpecial(-4122)

# Scroll the view down to roughly where the new row lives and move the
# active selection to the next empty detail cell, same as the author left
# the sheet after typing the new row in Excel.
$win = $excel.ActiveWindow
$win.ScrollRow = 28
$win.ScrollColumn = 2
$ws.Range("E43").Select()
